$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.871.47"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.625.27"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'211.04"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'23.42"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "'0.0878"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.855.93"
$ws.Range("D13").Value = "1.616.95"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "'65.32"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "27.857.38"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'229.48"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'7.64"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "0.0₃0721"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  -5.34%  "
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("D25").Value = "'155.03"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").Value = "'6.91"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'15.50"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("D34").Value = "1.392.48"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("E36").Value = "  +11.79%  "
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'0.858"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("D45").Value = "'65.67"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.765.30"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "'2.17"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "'87.86"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  +1.25%  "

foreach ($addr in @("D4","D5","D8","D11","D16","D18","D19","D25","D26","D28","D40","D45","D47","D48")) {
    $ws.Range($addr).Style = "Normal"
}
